{"js": "// Splits three single-run \"Programa\"/\"Bibliografia\" paragraphs into\n// multiple <w:t>/<w:br/> children of ONE run (matching the target OOXML),\n// one <w:br/> between each topic/reference line.\n\nconst RUN_INNER_PT = \"<w:t>Conceitos introdut\u00f3rios: composi\u00e7\u00e3o qu\u00edmica e mineral\u00f3gica do solo, horizontes de solo, unidades aqu\u00edferas</w:t><w:br/><w:t>Principais poluentes dos solos e \u00e1guas subterr\u00e2neas</w:t><w:br/><w:t>Legisla\u00e7\u00e3o: valores orientadores</w:t><w:br/><w:t xml:space=\\\"preserve\\\">Intera\u00e7\u00e3o entre solo e poluentes </w:t><w:br/><w:t>Fluxo de \u00e1gua em zonas n\u00e3o saturadas e saturadas</w:t><w:br/><w:t>Fluxo de \u00e1gua em fraturas de rocha</w:t><w:br/><w:t xml:space=\\\"preserve\\\">Transporte de subst\u00e2ncias misc\u00edveis na \u00e1gua no solo e \u00e1guas subterr\u00e2neas: mecanismos, modelos e solu\u00e7\u00f5es </w:t><w:br/><w:t xml:space=\\\"preserve\\\">Transporte de subst\u00e2ncias n\u00e3o misc\u00edveis na \u00e1gua (NAPL) no solo e \u00e1guas subterr\u00e2neas: mecanismos, modelos e solu\u00e7\u00f5es </w:t><w:br/><w:t xml:space=\\\"preserve\\\">Gerenciamento de \u00e1reas contaminadas: investiga\u00e7\u00e3o preliminar, t\u00e9cnicas de investiga\u00e7\u00e3o geol\u00f3gica-geot\u00e9cnica; modelo conceitual; an\u00e1lise de risco e t\u00e9cnicas de interven\u00e7\u00e3o </w:t><w:br/><w:t>Exemplos de aplica\u00e7\u00e3o em problemas geoambientais</w:t><w:br/><w:t>A disciplina pode contar com viagens did\u00e1ticas para complementa\u00e7\u00e3o do conte\u00fado da disciplina</w:t>\";\nconst RUN_INNER_EN = \"<w:rPr><w:i/></w:rPr><w:t>Introductory concepts: chemical and mineralogical composition of the soil, soil horizons, aquifer units</w:t><w:br/><w:t>Main pollutants of soil and groundwater</w:t><w:br/><w:t>Legislation and guiding values</w:t><w:br/><w:t>Interaction between soil and pollutants</w:t><w:br/><w:t>Water flow in unsaturated and saturated zones</w:t><w:br/><w:t>Water flow in rock fractures</w:t><w:br/><w:t>Transport of water-miscible substances in soil and groundwater: mechanisms, models and solutions</w:t><w:br/><w:t>Transport of non-aqueous phase liquid (NAPL) in soil and groundwater: mechanisms, models and solutions</w:t><w:br/><w:t>Management of contaminated areas: preliminary investigation, geological-geotechnical investigation techniques; conceptual model; risk analysis and intervention techniques</w:t><w:br/><w:t>Examples of application to geoenvironmental problems</w:t><w:br/><w:t>The discipline may have didactic trips to complement the content of the discipline</w:t>\";\nconst RUN_INNER_BIB = \"<w:t>BOSCOV, M.E.G. Geotecnia Ambiental. Oficina de Textos, 2008. 248 p.</w:t><w:br/><w:t>CETESB. Decis\u00e3o da Diretoria N. 38/2017/C, de 7 de fevereiro de 2017. 65 p.</w:t><w:br/><w:t>OLIVEIRA, A. M. S.; JER\u00d4NIMO, J. Geologia de Engenharia e Ambiental, ABGE, 2018. 912 p.</w:t><w:br/><w:t xml:space=\\\"preserve\\\">SHARMA, H. D.; REDDY, K. R. Geoenvironmental engineering, Wiley, 2004. 992p. </w:t><w:br/><w:t>YONG, R. N. Geoenvironmental engineering. contaminated soils, Pollutant fate and migration. CRC Press, 2001. 307p.</w:t>\";\n\nconst PREFIX_PT = \"Conceitos introdut\u00f3rios\";\nconst PREFIX_EN = \"Introductory concepts\";\nconst PREFIX_BIB = \"BOSCOV, M.E.G.\";\n\nfunction wrapOoxml(runInnerXml) {\n  return '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<?mso-application progid=\"Word.Document\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    '<pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p><w:r>' + runInnerXml + '</w:r></w:p>' +\n    '<w:sectPr/>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  { prefix: PREFIX_PT, inner: RUN_INNER_PT },\n  { prefix: PREFIX_EN, inner: RUN_INNER_EN },\n  { prefix: PREFIX_BIB, inner: RUN_INNER_BIB },\n];\n\nfor (const { prefix, inner } of targets) {\n  const para = paragraphs.items.find((p) => p.text.indexOf(prefix) === 0);\n  if (!para) {\n    throw new Error(\"Could not locate paragraph starting with: \" + prefix);\n  }\n  const range = para.getRange();\n  range.insertOoxml(wrapOoxml(inner), Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Splits three single-run paragraphs (\"Programa\" PT/EN and \"Bibliografia\")\n# into one run per paragraph containing multiple <w:t> runs of text joined\n# by line breaks (Chr(11) / vbVerticalTab -> <w:br/> in OOXML), one break\n# between each topic / reference line.\n\n$d = $word.ActiveDocument\n\n$lineBreak = [char]11\n\n$targets = @(\n    @{\n        Prefix   = 'Conceitos introdut\u00f3rios'\n        Segments = @(\n        'Conceitos introdut\u00f3rios: composi\u00e7\u00e3o qu\u00edmica e mineral\u00f3gica do solo, horizontes de solo, unidades aqu\u00edferas',\n        'Principais poluentes dos solos e \u00e1guas subterr\u00e2neas',\n        'Legisla\u00e7\u00e3o: valores orientadores',\n        'Intera\u00e7\u00e3o entre solo e poluentes ',\n        'Fluxo de \u00e1gua em zonas n\u00e3o saturadas e saturadas',\n        'Fluxo de \u00e1gua em fraturas de rocha',\n        'Transporte de subst\u00e2ncias misc\u00edveis na \u00e1gua no solo e \u00e1guas subterr\u00e2neas: mecanismos, modelos e solu\u00e7\u00f5es ',\n        'Transporte de subst\u00e2ncias n\u00e3o misc\u00edveis na \u00e1gua (NAPL) no solo e \u00e1guas subterr\u00e2neas: mecanismos, modelos e solu\u00e7\u00f5es ',\n        'Gerenciamento de \u00e1reas contaminadas: investiga\u00e7\u00e3o preliminar, t\u00e9cnicas de investiga\u00e7\u00e3o geol\u00f3gica-geot\u00e9cnica; modelo conceitual; an\u00e1lise de risco e t\u00e9cnicas de interven\u00e7\u00e3o ',\n        'Exemplos de aplica\u00e7\u00e3o em problemas geoambientais',\n        'A disciplina pode contar com viagens did\u00e1ticas para complementa\u00e7\u00e3o do conte\u00fado da disciplina'\n    )\n    },\n    @{\n        Prefix   = 'Introductory concepts'\n        Segments = @(\n        'Introductory concepts: chemical and mineralogical composition of the soil, soil horizons, aquifer units',\n        'Main pollutants of soil and groundwater',\n        'Legislation and guiding values',\n        'Interaction between soil and pollutants',\n        'Water flow in unsaturated and saturated zones',\n        'Water flow in rock fractures',\n        'Transport of water-miscible substances in soil and groundwater: mechanisms, models and solutions',\n        'Transport of non-aqueous phase liquid (NAPL) in soil and groundwater: mechanisms, models and solutions',\n        'Management of contaminated areas: preliminary investigation, geological-geotechnical investigation techniques; conceptual model; risk analysis and intervention techniques',\n        'Examples of application to geoenvironmental problems',\n        'The discipline may have didactic trips to complement the content of the discipline'\n    )\n    },\n    @{\n        Prefix   = 'BOSCOV, M.E.G.'\n        Segments = @(\n        'BOSCOV, M.E.G. Geotecnia Ambiental. Oficina de Textos, 2008. 248 p.',\n        'CETESB. Decis\u00e3o da Diretoria N. 38/2017/C, de 7 de fevereiro de 2017. 65 p.',\n        'OLIVEIRA, A. M. S.; JER\u00d4NIMO, J. Geologia de Engenharia e Ambiental, ABGE, 2018. 912 p.',\n        'SHARMA, H. D.; REDDY, K. R. Geoenvironmental engineering, Wiley, 2004. 992p. ',\n        'YONG, R. N. Geoenvironmental engineering. contaminated soils, Pollutant fate and migration. CRC Press, 2001. 307p.'\n    )\n    }\n)\n\nforeach ($target in $targets) {\n    $newText = [string]::Join($lineBreak, $target.Segments)\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.StartsWith($target.Prefix)) {\n            $p.Range.Text = $newText\n            break\n        }\n    }\n}\n"}
